$d = $word.ActiveDocument

# --- Title / TitleChar: drop the condensed-spacing / kerning direct
#     formatting that used to sit on the big display heading. ---
$title = $d.Styles("Title")
$title.Font.Spacing = 0
$title.Font.Kerning = 0

$titleChar = $d.Styles("TitleChar")
$titleChar.Font.Spacing = 0
$titleChar.Font.Kerning = 0

# --- Author: now inherits from Title instead of duplicating its own
#     centering; gets an explicit, smaller run size instead. ---
$author = $d.Styles("Author")
$author.BaseStyle = $title
$author.ParagraphFormat.Alignment = 0
$author.Font.Size = 12
$author.Font.SizeBi = 12

# --- Date: same treatment as Author. ---
$date = $d.Styles("Date")
$date.BaseStyle = $title
$date.ParagraphFormat.Alignment = 0
$date.Font.Size = 12
$date.Font.SizeBi = 12
